$d = $word.ActiveDocument

# 1. Merge the "lithographic stones of head" / bookmark / "ing 8442;" runs into a
#    single contiguous run, which also removes the old _GoBack bookmark that sat
#    between them.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("lithographic stones of heading 8442;", $true, $false, $false, $false, $false, $true, 1, $false, "lithographic stones of heading 8442;", 2) | Out-Null

# 2. Re-create the _GoBack bookmark collapsed at the very start of the document,
#    right before the "Section XIII" heading run.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Inserting a throwaway character at the very start first (then bookmarking at
# position 1, then deleting the throwaway character) avoids an engine quirk
# where bookmarking a zero-length range exactly at offset 0 mis-places the
# bookmarkEnd element elsewhere in the document.
$head = $d.Range(0, 0)
$head.InsertBefore("x")
$anchor = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $anchor)
$d.Range(0, 1).Delete()

# 3. "Normal in Table" style: reduce the run font size from 9pt to 8pt.
$tableStyle = $d.Styles.Item("Normal in Table")
$tableStyle.Font.Size = 8
